$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 15 and 16 (last two data rows no longer present in target)
$ws.Rows.Item(15).Delete() | Out-Null
$ws.Rows.Item(15).Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'life-dev/main'
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = 'impression'
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = 'channel, page_url, os_name, impression_type'
$ws.Range("G2").Value = 'Rround, https://life-dev.hectoinnovation.co.kr/main, iOS, 검색 창'
$ws.Range("H2").Value = 4

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 'life-dev/main'
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = 'impression'
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = 'channel, page_url, banner_text, banner_position, os_name, impression_type'
$ws.Range("G3").Value = 'Rround, https://life-dev.hectoinnovation.co.kr/main, ______포인트 12,500원 놓치고 있어요!___전국 날씨특파원, 오늘 날씨는?___여름맞이 체력 증진! 오운완 챌린지___, 라이프 메인 상단 카드 배너, iOS, 메인 상단 카드 배너'
$ws.Range("H3").Value = 6

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'life-dev/main'
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = 'pageview'
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = 'channel, page_url, os_name'
$ws.Range("G4").Value = 'Rround, https://life-dev.hectoinnovation.co.kr/main, iOS'
$ws.Range("H4").Value = 3

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'life-dev/main'
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = 'swipe'
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 'channel, page_url, swipe_area, swipe_direct, os_name'
$ws.Range("G5").Value = 'Rround, https://life-dev.hectoinnovation.co.kr/main, 메인 상단 카드 배너, left, iOS'
$ws.Range("H5").Value = 5

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 'life-dev/main'
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = 'swipe'
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = 'channel, page_url, swipe_area, swipe_direct, os_name'
$ws.Range("G6").Value = 'Rround, https://life-dev.hectoinnovation.co.kr/main, 메인 상단 카드 배너, left, iOS'
$ws.Range("H6").Value = 5

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 'life-dev/main'
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = 'swipe'
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = 'channel, page_url, swipe_area, swipe_direct, os_name'
$ws.Range("G7").Value = 'Rround, https://life-dev.hectoinnovation.co.kr/main, 메인 상단 카드 배너, left, iOS'
$ws.Range("H7").Value = 5

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'life-dev/main'
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 'swipe'
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = 'channel, page_url, swipe_area, swipe_direct, os_name'
$ws.Range("G8").Value = 'Rround, https://life-dev.hectoinnovation.co.kr/main, 메인 상단 카드 배너, left, iOS'
$ws.Range("H8").Value = 5

# Row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 'ecommerce-dev/product/detail/800'
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = 'click'
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = 'channel, page_url, tab_name, prd_code, prd_name, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, os_name'
$ws.Range("G9").Value = 'Rround, https://ecommerce-dev.hectoinnovation.co.kr/product/detail/800, 상품상세
, 800, 여성용 스킨핏 50수 투톤 모달 팬티 5P SET, 20,000원, 20,000원, 10%, 0, 0, #여성팬티___#50수팬티___#숙녀팬티___#여자팬티___#팬티세트___#모달팬티___#투톤팬티___#팬티, iOS'
$ws.Range("H9").Value = 12

# Row 10
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 'life-dev/main'
$ws.Range("C10").Value = '상품 찜하기'
$ws.Range("D10").Value = 'click'
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = 'channel, page_url, prd_code, prd_name, prd_brand, prd_price_final, prd_is_ad, os_name'
$ws.Range("G10").Value = 'Rround, https://life-dev.hectoinnovation.co.kr/main, 800, 여성용 스킨핏 50수 투톤 모달 팬티 5P SET, 마이그스토어, 20,000원, F, iOS'
$ws.Range("H10").Value = 8

# Row 11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 'life-dev/main'
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = 'click'
$ws.Range("E11").Value = '상품 더보기'
$ws.Range("F11").Value = 'channel, page_url, click_text, module_id, module_order, module_name, os_name'
$ws.Range("G11").Value = 'Rround, https://life-dev.hectoinnovation.co.kr/main, 상품 더보기, C-3, 33, commerce-category-ranking, iOS'
$ws.Range("H11").Value = 7

# Row 12
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 'ecommerce-dev/category/detail/543'
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = 'pageview'
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = 'channel, page_url, ctgr_id, os_name'
$ws.Range("G12").Value = 'Rround, https://ecommerce-dev.hectoinnovation.co.kr/category/detail/543, 543, iOS'
$ws.Range("H12").Value = 4

# Row 13
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 'life-dev/main'
$ws.Range("C13").Value = '상품'
$ws.Range("D13").Value = 'click'
$ws.Range("E13").Value = '드시모네 베이비스텝2 100억 생유산균 2박스'
$ws.Range("F13").Value = 'channel, page_url, click_text, module_id, module_order, prd_order, prd_code, prd_name, prd_brand, prd_price_final, prd_is_ad, el_order, module_name, os_name'
$ws.Range("G13").Value = 'Rround, https://life-dev.hectoinnovation.co.kr/main, 드시모네 베이비스텝2 100억 생유산균 2박스, C-3, 33, 2, 1030, 드시모네 베이비스텝2 100억 생유산균 2박스, 마이그스토어, 96,000원, F, 2, commerce-category-ranking, iOS'
$ws.Range("H13").Value = 14

# Row 14
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 'ecommerce-dev/product/detail/1030'
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = 'pageview'
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = 'channel, page_url, prd_code, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, os_name'
$ws.Range("G14").Value = 'Rround, https://ecommerce-dev.hectoinnovation.co.kr/product/detail/1030, 1030, 96,000원, 86,400원, 10%, 0, 0, #프로바이오틱스___#식품___#영양제___#드시모네___#베이비스텝___#박스___#생유산균___#건강식품, iOS'
$ws.Range("H14").Value = 10
